$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "nom" column (CE, column 83).
# This shifts the old "nom" (CE) -> CF and old "url_produit" (CF) -> CG,
# matching the dimension change from A1:CF206 to A1:CG206.
$ws.Columns.Item(83).Insert()

# New column header: latest price-check timestamp.
$ws.Range("CE1").Value = "2026-01-31 11:11:58"

# For the rows that already have numeric price history (rows 2-80),
# the newest check repeats the last known price (same as column CD).
# Copy the values across so the numeric type/precision is preserved exactly.
$src = $ws.Range("CD2:CD80")
$dst = $ws.Range("CE2:CE80")
$src.Copy($dst)
